$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.653.67"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.279.62"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'504.53"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "'128.66"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "2.297.29"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").Value = "2.685.94"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "54.730.72"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "2.283.90"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'307.27"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "'6.44"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").Value = "'0.994"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'7.45"
$ws.Range("D28").Value = "'171.30"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0701"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.63"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D34").Value = "'17.95"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'3.78"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.39"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'4.91"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "'127.07"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "'250.88"
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  +0.46%  "
